$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update DAMSLTag (column I) and DialogAct (column J) values for the rows
# re-annotated by the SGNN re-run following transcript clean-up.
$updates = @(
    @{ Row = 6; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 7; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 15; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 18; I = "sv"; J = "Statement-opinion" }
    @{ Row = 19; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 26; I = "aa"; J = "Agree/Accept" }
    @{ Row = 28; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 41; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 49; I = "sv"; J = "Statement-opinion" }
    @{ Row = 51; I = "b"; J = "Acknowledge (Backchannel)" }
    @{ Row = 58; I = "aa"; J = "Agree/Accept" }
    @{ Row = 64; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 83; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 88; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 93; I = "sv"; J = "Statement-opinion" }
    @{ Row = 98; I = "sv"; J = "Statement-opinion" }
    @{ Row = 115; I = "aa"; J = "Agree/Accept" }
    @{ Row = 121; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 123; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 142; I = "aa"; J = "Agree/Accept" }
    @{ Row = 156; I = "sv"; J = "Statement-opinion" }
    @{ Row = 157; I = "sv"; J = "Statement-opinion" }
    @{ Row = 159; I = "sv"; J = "Statement-opinion" }
    @{ Row = 160; I = "%"; J = "Uninterpretable" }
    @{ Row = 168; I = "b"; J = "Acknowledge (Backchannel)" }
    @{ Row = 175; I = "ba"; J = "Appreciation" }
    @{ Row = 181; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 182; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 200; I = "sv"; J = "Statement-opinion" }
    @{ Row = 221; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 222; I = "sv"; J = "Statement-opinion" }
    @{ Row = 235; I = "aa"; J = "Agree/Accept" }
    @{ Row = 239; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 240; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 273; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 284; I = "b"; J = "Acknowledge (Backchannel)" }
    @{ Row = 302; I = "sv"; J = "Statement-opinion" }
    @{ Row = 304; I = "sv"; J = "Statement-opinion" }
    @{ Row = 306; I = "sv"; J = "Statement-opinion" }
    @{ Row = 310; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 313; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 329; I = "sv"; J = "Statement-opinion" }
    @{ Row = 334; I = "sv"; J = "Statement-opinion" }
    @{ Row = 342; I = "%"; J = "Uninterpretable" }
    @{ Row = 351; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 370; I = "%"; J = "Uninterpretable" }
    @{ Row = 372; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 417; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 429; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 440; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 441; I = "ba"; J = "Appreciation" }
    @{ Row = 443; I = "sv"; J = "Statement-opinion" }
    @{ Row = 464; I = "b"; J = "Acknowledge (Backchannel)" }
    @{ Row = 467; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 469; I = "sv"; J = "Statement-opinion" }
    @{ Row = 479; I = "sv"; J = "Statement-opinion" }
    @{ Row = 492; I = "sv"; J = "Statement-opinion" }
    @{ Row = 505; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 515; I = "ba"; J = "Appreciation" }
    @{ Row = 528; I = "b"; J = "Acknowledge (Backchannel)" }
    @{ Row = 537; I = "b"; J = "Acknowledge (Backchannel)" }
    @{ Row = 547; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 563; I = "aa"; J = "Agree/Accept" }
    @{ Row = 567; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 569; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 571; I = "ba"; J = "Appreciation" }
    @{ Row = 584; I = "sv"; J = "Statement-opinion" }
    @{ Row = 588; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 589; I = "sv"; J = "Statement-opinion" }
    @{ Row = 596; I = "aa"; J = "Agree/Accept" }
    @{ Row = 598; I = "aa"; J = "Agree/Accept" }
    @{ Row = 606; I = "sv"; J = "Statement-opinion" }
    @{ Row = 609; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 613; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 626; I = "sv"; J = "Statement-opinion" }
    @{ Row = 644; I = "sv"; J = "Statement-opinion" }
    @{ Row = 652; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 658; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 662; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 665; I = "sv"; J = "Statement-opinion" }
    @{ Row = 693; I = "ba"; J = "Appreciation" }
    @{ Row = 706; I = "b"; J = "Acknowledge (Backchannel)" }
    @{ Row = 709; I = "aa"; J = "Agree/Accept" }
    @{ Row = 713; I = "sv"; J = "Statement-opinion" }
    @{ Row = 715; I = "b"; J = "Acknowledge (Backchannel)" }
    @{ Row = 717; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 718; I = "sd"; J = "Statement-non-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
